# Update "想去人数" (want-to-go count) column F on the 展览 and 全部类型
# sheets to the refreshed scrape values.

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item("展览")
$sheetAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1) — F column updates keyed by row number
$exhibitUpdates = @{
    2  = 188
    3  = 5343
    6  = 23
    7  = 604
    8  = 573
    9  = 1051
    11 = 1471
    12 = 4280
    14 = 191
    17 = 3412
    18 = 168
    19 = 1089
    20 = 103
    22 = 199
    23 = 125
    24 = 40
    25 = 139
    27 = 308
    29 = 56
    30 = 16
    31 = 28
}

foreach ($row in $exhibitUpdates.Keys) {
    $sheetExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# 全部类型 (sheet4) — F column updates keyed by row number
$allUpdates = @{
    2  = 188
    4  = 5343
    7  = 23
    8  = 604
    9  = 573
    10 = 1051
    12 = 1471
    13 = 4280
    15 = 191
    18 = 3412
    19 = 168
    20 = 1089
    21 = 103
    23 = 199
    24 = 125
    25 = 40
    26 = 139
    28 = 308
    30 = 56
    31 = 16
    32 = 28
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
